# Add a "Career_Profiling_Flag" column (T) to the "August" sheet.
# The flag is 1 for any row whose "Person tag" (column E) text contains
# the substring "Career Profiling Engaged", and 0 otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("August")

# Determine the last used data row on the sheet.
$lastRow = $ws.UsedRange.Rows.Count

# Header for the new column.
$ws.Cells.Item(1, 20).Value2 = "Career_Profiling_Flag"

# Walk every data row (2 .. lastRow, inclusive of the trailing totals row)
# and flag whether the Person tag column mentions Career Profiling engagement.
for ($r = 2; $r -le $lastRow; $r++) {
    $tag = $ws.Cells.Item($r, 5).Value2
    $flag = 0
    if ($tag -ne $null -and $tag.ToString().Contains("Career Profiling Engaged")) {
        $flag = 1
    }
    $ws.Cells.Item($r, 20).Value2 = $flag
}
